$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1346003333333333
$ws.Range("H2").Value = 0.403801
$ws.Range("I2").Value = 0.009651054304565105
$ws.Range("J2").Value = 0.009651054304565105
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 0.05770895071433334
$ws.Range("R2").Value = 0.5193805564290001
$ws.Range("S2").Value = [double]"3.976720737041929E-05"
$ws.Range("T2").Value = [double]"3.976720737041929E-05"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1346003333333333
$ws.Range("H3").Value = 0.403801
$ws.Range("I3").Value = 0.009651054304565105
$ws.Range("J3").Value = 0.009651054304565105
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 10.798478242279
$ws.Range("R3").Value = 97.18630418051099
$ws.Range("S3").Value = 0.007441225637100556
$ws.Range("T3").Value = 0.007441225637100556
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1346003333333333
$ws.Range("H4").Value = 0.403801
$ws.Range("I4").Value = 0.009651054304565105
$ws.Range("J4").Value = 0.009651054304565105
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 3.149126582642556
$ws.Range("R4").Value = 28.342139243783
$ws.Range("S4").Value = 0.002170061460094129
$ws.Range("T4").Value = 0.002170061460094129
$ws.Range("I5").Value = 0.8124788779145131
$ws.Range("J5").Value = 0.8124788779145132
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 4.858257143971001
$ws.Range("R5").Value = 43.72431429573901
$ws.Range("S5").Value = 0.003347822424626588
$ws.Range("T5").Value = 0.003347822424626589
$ws.Range("I6").Value = 0.8124788779145131
$ws.Range("J6").Value = 0.8124788779145132
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("Q6").Value = 909.0753412630888
$ws.Range("R6").Value = 8181.6780713678
$ws.Range("S6").Value = 0.6264433361524439
$ws.Range("T6").Value = 0.626443336152444
$ws.Range("I7").Value = 0.8124788779145131
$ws.Range("J7").Value = 0.8124788779145132
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 265.1108108536837
$ws.Range("R7").Value = 2385.997297683153
$ws.Range("S7").Value = 0.1826877193374427
$ws.Range("T7").Value = 0.1826877193374427
$ws.Range("G8").Value = 2.4807
$ws.Range("H8").Value = 7.4421
$ws.Range("I8").Value = 0.1778700677809217
$ws.Range("J8").Value = 0.1778700677809217
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 1.0635827601
$ws.Range("R8").Value = 9.5722448409
$ws.Range("S8").Value = 0.0007329143166346726
$ws.Range("T8").Value = 0.0007329143166346726
$ws.Range("G9").Value = 2.4807
$ws.Range("H9").Value = 7.4421
$ws.Range("I9").Value = 0.1778700677809217
$ws.Range("J9").Value = 0.1778700677809217
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 199.0172261259
$ws.Range("R9").Value = 1791.1550351331
$ws.Range("S9").Value = 0.1371426651094625
$ws.Range("T9").Value = 0.1371426651094625
$ws.Range("G10").Value = 2.4807
$ws.Range("H10").Value = 7.4421
$ws.Range("I10").Value = 0.1778700677809217
$ws.Range("J10").Value = 0.1778700677809217
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 58.03877390270001
$ws.Range("R10").Value = 522.3489651243
$ws.Range("S10").Value = 0.03999448835482458
$ws.Range("T10").Value = 0.03999448835482458
